$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8241867217947743
$ws.Range("C2").Value = 0.1998266203563048
$ws.Range("E2").Value = 0.1136866213164893
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.5498833845211664
$ws.Range("H2").Value = 0.6254872185886313
$ws.Range("M2").Value = 0.3424943153330489
$ws.Range("B3").Value = 0.7251062635903622
$ws.Range("C3").Value = 0.1760536258992431
$ws.Range("E3").Value = 0.1081556909239438
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.5359160217344936
$ws.Range("H3").Value = 0.6255920788342166
$ws.Range("M3").Value = 0.305373343264236
$ws.Range("B4").Value = 0.6643145236133989
$ws.Range("C4").Value = 0.1613864964857328
$ws.Range("E4").Value = 0.1048574224325662
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.52804632916191
$ws.Range("H4").Value = 0.6262137786452655
$ws.Range("M4").Value = 0.28268572673813
$ws.Range("B5").Value = 0.6395524364540961
$ws.Range("C5").Value = 0.1553918225393431
$ws.Range("E5").Value = 0.1035376173745632
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.5250149041233954
$ws.Range("H5").Value = 0.6266064251323371
$ws.Range("M5").Value = 0.2734662246114112
$ws.Range("B6").Value = 0.6354413771467193
$ws.Range("C6").Value = 0.1543953428686109
$ws.Range("E6").Value = 0.1033199212597182
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.5245220797697243
$ws.Range("H6").Value = 0.6266800124445098
$ws.Range("M6").Value = 0.2719368832624767
$ws.Range("B7").Value = 0.663980528770054
$ws.Range("C7").Value = 0.1613057218643235
$ws.Range("E7").Value = 0.1048395252383116
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.5280047383007513
$ws.Range("H7").Value = 0.6262185111219623
$ws.Range("M7").Value = 0.2825612851066523
$ws.Range("B8").Value = 0.7900145717982809
$ws.Range("C8").Value = 0.191644301228024
$ws.Range("E8").Value = 0.1117590789784515
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.5449195876979758
$ws.Range("H8").Value = 0.625407194435553
$ws.Range("M8").Value = 0.3296729537806797
$ws.Range("B9").Value = 1.037532426757934
$ws.Range("C9").Value = 0.2505823831147325
$ws.Range("E9").Value = 0.1261187909381078
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.5837894590630697
$ws.Range("H9").Value = 0.6282770915465647
$ws.Range("M9").Value = 0.4229167418350102
$ws.Range("B10").Value = 1.219648908285251
$ws.Range("C10").Value = 0.2935540916710124
$ws.Range("E10").Value = 0.1371726235063448
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.6159571862790756
$ws.Range("H10").Value = 0.6331600600547347
$ws.Range("M10").Value = 0.4919898512155072
$ws.Range("B11").Value = 1.30256680139621
$ws.Range("C11").Value = 0.3130336421376683
$ws.Range("E11").Value = 0.1423155000526748
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.6314043384397792
$ws.Range("H11").Value = 0.6359958531092218
$ws.Range("M11").Value = 0.5235461676762583
$ws.Range("B12").Value = 1.333976540296533
$ws.Range("C12").Value = 0.3204002944629565
$ws.Range("E12").Value = 0.1442798002885368
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.6373731307708681
$ws.Range("H12").Value = 0.6371590139954151
$ws.Range("M12").Value = 0.5355157631142191
$ws.Range("B13").Value = 1.327211419844843
$ws.Range("C13").Value = 0.318814190909336
$ws.Range("E13").Value = 0.1438560001039306
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.6360823030601921
$ws.Range("H13").Value = 0.6369045194509226
$ws.Range("M13").Value = 0.5329370025308435
$ws.Range("B14").Value = 1.305150687452112
$ws.Range("C14").Value = 0.3136398979400212
$ws.Range("E14").Value = 0.142476765414898
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.6318929908645998
$ws.Range("H14").Value = 0.6360897517846809
$ws.Range("M14").Value = 0.5245305114288925
$ws.Range("B15").Value = 1.291639228541897
$ws.Range("C15").Value = 0.3104692145566901
$ws.Range("E15").Value = 0.141634143210986
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.6293425184433943
$ws.Range("H15").Value = 0.6356023415935397
$ws.Range("M15").Value = 0.5193839013291495
$ws.Range("B16").Value = 1.214231493830766
$ws.Range("C16").Value = 0.2922796803521521
$ws.Range("E16").Value = 0.1368388568797911
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.6149642486729761
$ws.Range("H16").Value = 0.6329871801982847
$ws.Range("M16").Value = 0.4899303268477979
$ws.Range("B17").Value = 1.166762943156016
$ws.Range("C17").Value = 0.2811034576815246
$ws.Range("E17").Value = 0.1339266603622633
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.6063537046602221
$ws.Range("H17").Value = 0.6315409222534925
$ws.Range("M17").Value = 0.4718963747417035
$ws.Range("B18").Value = 1.139467072265404
$ws.Range("C18").Value = 0.2746687412250992
$ws.Range("E18").Value = 0.1322623918047938
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.6014776196077776
$ws.Range("H18").Value = 0.6307668565822979
$ws.Range("M18").Value = 0.4615363593519817
$ws.Range("B19").Value = 1.130226320478641
$ws.Range("C19").Value = 0.2724889487816426
$ws.Range("E19").Value = 0.1317007347817096
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.599839726582502
$ws.Range("H19").Value = 0.6305146669944577
$ws.Range("M19").Value = 0.4580307901299108
$ws.Range("B20").Value = 1.171815348291204
$ws.Range("C20").Value = 0.2822938529699002
$ws.Range("E20").Value = 0.1342355535629309
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.6072623809002664
$ws.Range("H20").Value = 0.6316888917013728
$ws.Range("M20").Value = 0.4738148079359945
$ws.Range("B21").Value = 1.311630172103492
$ws.Range("C21").Value = 0.3151599798203222
$ws.Range("E21").Value = 0.1428814212061624
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.6331202378749765
$ws.Range("H21").Value = 0.6363266374823127
$ws.Range("M21").Value = 0.5269991570538082
$ws.Range("B22").Value = 1.403068725980802
$ws.Range("C22").Value = 0.3365826078656085
$ws.Range("E22").Value = 0.1486301096337641
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.6507162138709646
$ws.Range("H22").Value = 0.6398786661339102
$ws.Range("M22").Value = 0.5618744782308767
$ws.Range("B23").Value = 1.354260564600736
$ws.Range("C23").Value = 0.3251541812781227
$ws.Range("E23").Value = 0.1455528326649826
$ws.Range("F23").Value = 0.7472568307915566
$ws.Range("G23").Value = 0.6412604452479229
$ws.Range("H23").Value = 0.6379348980609052
$ws.Range("M23").Value = 0.543250039858421
$ws.Range("B24").Value = 1.16953117329308
$ws.Range("C24").Value = 0.2817557044122339
$ws.Range("E24").Value = 0.1340958718517982
$ws.Range("F24").Value = 0.6416283278902313
$ws.Range("G24").Value = 0.6068513373240592
$ws.Range("H24").Value = 0.6316218159643086
$ws.Range("M24").Value = 0.4729474595923477
$ws.Range("B25").Value = 0.9705289090056795
$ws.Range("C25").Value = 0.2346967009590628
$ws.Range("E25").Value = 0.1221472341768859
$ws.Range("F25").Value = 0.5279251897347308
$ws.Range("G25").Value = 0.5726499605915194
$ws.Range("H25").Value = 0.6270174702410714
$ws.Range("M25").Value = 0.3975954900624856
